## "Audit info only in parent row"
##
## Column A ("Match ID") groups several rows together (one match can span
## several ledger rows). Column B ("Audit Info" - Match Type / L/C / Lender
## & Borrower amount) used to be repeated on every row of the group. This
## trims it down so the Audit Info text (and its wrap-text style) is kept
## only on the first ("parent") row of each Match ID group, and removed
## from every other row that belongs to the same group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$prevMatchId = ""
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $matchId = $ws.Cells.Item($r, 1).Text

    if ($matchId -ne "") {
        if ($matchId -eq $prevMatchId) {
            # Same Match ID as the row above -> this is a child row, not the
            # parent row of the group, so its Audit Info is cleared (value +
            # formatting, matching how the parent-row copy looks untouched).
            $ws.Cells.Item($r, 2).Clear()
        }
        $prevMatchId = $matchId
    }
}
